$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b4Query = @'
MATCH (f:file)-->(parent)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (diag:diagnosis)-->(c)
WHERE demo.breed IN ['American Staffordshire Terrier']
MATCH (s:study)<--(c)
WHERE s.clinical_study_designation IN ['COTC022']
OPTIONAL MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@

$b5Query = @'
MATCH (f:file)-->(s:study)
WHERE s.clinical_study_designation IN ['COTC022']
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$ws.Range("B4").Value = $b4Query
$ws.Range("B5").Value = $b5Query

# Update the saved view state (scroll position, zoom, active selection)
$ws.Activate()
$excel.Goto($ws.Range("A5"), $true)
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ZoomScaleNormal = 70
$ws.Range("B6").Select()

# Theme was refreshed to the newer "Office 2013 - 2022" palette/name on re-save
try {
    $wb.Theme.Name = "Office 2013 - 2022 Theme"
    $wb.Theme.ThemeColorScheme.Name = "Office 2013 - 2022"
    $wb.Theme.ThemeFontScheme.Name = "Office 2013 - 2022"
} catch {
}
